# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of an existing header cell (bold, centered, bordered)
# onto the three new header cells so they match the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (rows 2-62): season record values for every player row ---
$lastRow = 62
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 65
    $ws.Cells.Item($r, 31).Value = 97
    $ws.Cells.Item($r, 32).Value = 0
}
